{"js": "// Apply a yellow highlight to the requirement line about adding a new\n// product, matching the XML diff's addition of <w:highlight w:val=\"yellow\"/>\n// to that run's rPr.\nconst targetText =\n  \"El sistema debe agregar un nuevo producto. (panel principal)\";\n\nconst results = context.document.body.search(targetText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found: \" + targetText);\n}\n\n// Highlight every match (expected to be exactly one) in yellow.\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].font.highlightColor = \"Yellow\";\n}\n\nawait context.sync();\n", "ps1": "# Apply a yellow highlight to the requirement line about adding a new\n# product, matching the XML diff's addition of <w:highlight w:val=\"yellow\"/>\n# to that run's rPr.\n$d = $word.ActiveDocument\n\n$targetText = \"El sistema debe agregar un nuevo producto. (panel principal)\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $targetText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n\nwhile ($find.Execute()) {\n    $rng.HighlightColorIndex = 7  # wdYellow\n    $rng.Collapse(0)  # wdCollapseEnd, so Execute continues past this match\n}\n"}
